$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("B7").Value = 5.188199999999999
$ws.Range("A9").Value = -21.9494
$ws.Range("B12").Value = 5.610699999999992
$ws.Range("A18").Value = -22.26880000000001
$ws.Range("A20").Value = -19.37249999999998
$ws.Range("B26").Value = 4.222200000000003
$ws.Range("A27").Value = -21.65729999999997
$ws.Range("B27").Value = 4.940900000000002
$ws.Range("B29").Value = 5.1858
$ws.Range("B37").Value = 8.824500000000002
$ws.Range("B38").Value = 4.479100000000002
$ws.Range("B51").Value = 5.647699999999999
$ws.Range("B55").Value = 5.624299999999997
$ws.Range("A69").Value = -21.73259999999999
$ws.Range("B69").Value = 5.683199999999995
$ws.Range("B70").Value = 6.4555
$ws.Range("A76").Value = -19.56139999999998
$ws.Range("A82").Value = -22.00249999999999
$ws.Range("B83").Value = 5.438399999999997
$ws.Range("B102").Value = 8.552300000000004
